$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.187710165977478
$ws.Range("B1").Value = 2.203523635864258
$ws.Range("C1").Value = 6.437413215637207
$ws.Range("D1").Value = 2.304277658462524
$ws.Range("E1").Value = 1.192286729812622
